$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 12: Miles McElhany
$ws.Cells.Item(12, 2).Value2 = 343
$ws.Cells.Item(12, 3).Value2 = "Miles McElhany"
$ws.Cells.Item(12, 4).Value2 = 14
$ws.Cells.Item(12, 5).Value2 = 42
$ws.Cells.Item(12, 6).Value2 = $null
$ws.Cells.Item(12, 7).Value2 = 17
$ws.Cells.Item(12, 8).Value2 = 71
$ws.Cells.Item(12, 9).Value2 = 43
$ws.Cells.Item(12, 10).Value2 = $null
$ws.Cells.Item(12, 11).Value2 = $null
$ws.Cells.Item(12, 12).Value2 = $null
$ws.Cells.Item(12, 13).Value2 = $null
$ws.Cells.Item(12, 14).Value2 = $null
$ws.Cells.Item(12, 15).Value2 = $null
$ws.Cells.Item(12, 16).Value2 = 147
$ws.Cells.Item(12, 17).Value2 = $null
$ws.Cells.Item(12, 18).Value2 = $null
$ws.Cells.Item(12, 19).Value2 = $null
$ws.Cells.Item(12, 20).Formula = "=147+71"

# Row 13: Andy Halbert
$ws.Cells.Item(13, 2).Value2 = 107
$ws.Cells.Item(13, 3).Value2 = "Andy Halbert"
$ws.Cells.Item(13, 4).Value2 = $null
$ws.Cells.Item(13, 5).Value2 = $null
$ws.Cells.Item(13, 6).Value2 = $null
$ws.Cells.Item(13, 7).Value2 = $null
$ws.Cells.Item(13, 8).Value2 = $null
$ws.Cells.Item(13, 9).Value2 = $null
$ws.Cells.Item(13, 10).Value2 = $null
$ws.Cells.Item(13, 11).Value2 = $null
$ws.Cells.Item(13, 12).Value2 = 89.5
$ws.Cells.Item(13, 13).Value2 = $null
$ws.Cells.Item(13, 14).Value2 = 128
$ws.Cells.Item(13, 15).Value2 = $null
$ws.Cells.Item(13, 16).Value2 = $null
$ws.Cells.Item(13, 17).Value2 = $null
$ws.Cells.Item(13, 18).Value2 = $null
$ws.Cells.Item(13, 19).Value2 = $null
$ws.Cells.Item(13, 20).Formula = "=128+89.5"

# Row 14: Darrin Griffin
$ws.Cells.Item(14, 2).Value2 = 23
$ws.Cells.Item(14, 3).Value2 = "Darrin Griffin"
$ws.Cells.Item(14, 4).Value2 = $null
$ws.Cells.Item(14, 5).Value2 = $null
$ws.Cells.Item(14, 6).Value2 = $null
$ws.Cells.Item(14, 7).Value2 = $null
$ws.Cells.Item(14, 8).Value2 = $null
$ws.Cells.Item(14, 9).Value2 = $null
$ws.Cells.Item(14, 10).Value2 = $null
$ws.Cells.Item(14, 11).Value2 = $null
$ws.Cells.Item(14, 12).Value2 = $null
$ws.Cells.Item(14, 13).Value2 = $null
$ws.Cells.Item(14, 14).Value2 = $null
$ws.Cells.Item(14, 15).Value2 = $null
$ws.Cells.Item(14, 16).Value2 = $null
$ws.Cells.Item(14, 17).Value2 = 203.5
$ws.Cells.Item(14, 18).Value2 = $null
$ws.Cells.Item(14, 19).Value2 = $null
$ws.Cells.Item(14, 20).Formula = "=203.5"

# Row 15: David Pearce
$ws.Cells.Item(15, 2).Value2 = 17
$ws.Cells.Item(15, 3).Value2 = "David Pearce"
$ws.Cells.Item(15, 4).Value2 = $null
$ws.Cells.Item(15, 5).Value2 = $null
$ws.Cells.Item(15, 6).Value2 = $null
$ws.Cells.Item(15, 7).Value2 = $null
$ws.Cells.Item(15, 8).Value2 = $null
$ws.Cells.Item(15, 9).Value2 = $null
$ws.Cells.Item(15, 10).Value2 = $null
$ws.Cells.Item(15, 11).Value2 = $null
$ws.Cells.Item(15, 12).Value2 = 44.5
$ws.Cells.Item(15, 13).Value2 = $null
$ws.Cells.Item(15, 14).Value2 = 158.5
$ws.Cells.Item(15, 15).Value2 = $null
$ws.Cells.Item(15, 16).Value2 = $null
$ws.Cells.Item(15, 17).Value2 = $null
$ws.Cells.Item(15, 18).Value2 = $null
$ws.Cells.Item(15, 19).Value2 = $null
$ws.Cells.Item(15, 20).Formula = "=158.5+44.5"

# Row 16: Kevin Nanthrup
$ws.Cells.Item(16, 2).Value2 = 82
$ws.Cells.Item(16, 3).Value2 = "Kevin Nanthrup"
$ws.Cells.Item(16, 4).Value2 = 106
$ws.Cells.Item(16, 5).Value2 = $null
$ws.Cells.Item(16, 6).Value2 = 84
$ws.Cells.Item(16, 7).Value2 = 55
$ws.Cells.Item(16, 8).Value2 = $null
$ws.Cells.Item(16, 9).Value2 = $null
$ws.Cells.Item(16, 10).Value2 = $null
$ws.Cells.Item(16, 11).Value2 = $null
$ws.Cells.Item(16, 12).Value2 = $null
$ws.Cells.Item(16, 13).Value2 = $null
$ws.Cells.Item(16, 14).Value2 = $null
$ws.Cells.Item(16, 15).Value2 = $null
$ws.Cells.Item(16, 16).Value2 = $null
$ws.Cells.Item(16, 17).Value2 = $null
$ws.Cells.Item(16, 18).Value2 = $null
$ws.Cells.Item(16, 19).Value2 = $null
$ws.Cells.Item(16, 20).Formula = "=106+84"

# Row 17: Micah Kudo
$ws.Cells.Item(17, 2).Value2 = 84
$ws.Cells.Item(17, 3).Value2 = "Micah Kudo"
$ws.Cells.Item(17, 4).Value2 = $null
$ws.Cells.Item(17, 5).Value2 = 87
$ws.Cells.Item(17, 6).Value2 = $null
$ws.Cells.Item(17, 7).Value2 = $null
$ws.Cells.Item(17, 8).Value2 = 95
$ws.Cells.Item(17, 9).Value2 = 94
$ws.Cells.Item(17, 10).Value2 = $null
$ws.Cells.Item(17, 11).Value2 = $null
$ws.Cells.Item(17, 12).Value2 = $null
$ws.Cells.Item(17, 13).Value2 = $null
$ws.Cells.Item(17, 14).Value2 = $null
$ws.Cells.Item(17, 15).Value2 = $null
$ws.Cells.Item(17, 16).Value2 = $null
$ws.Cells.Item(17, 17).Value2 = $null
$ws.Cells.Item(17, 18).Value2 = $null
$ws.Cells.Item(17, 19).Value2 = $null
$ws.Cells.Item(17, 20).Formula = "=95+94"
